# Generate Report for Handoff
# Updates the "Ready for handoff" status (and related timestamps / error
# details) for the 75020932-... and d625f10b-... files across the
# Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$readyForHandoff = "Ready for handoff"

$msg75020932 = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e7f541279fc806ce39ecf1a53293d5c57c22ffd4/e2e/75020932-f74b-4555-a289-c249f94d18df.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/032a00e7b8b546b20e5ae53fa8d86a473852b56c/e2e/75020932-f74b-4555-a289-c249f94d18df.md."
$msgD625f10b = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e7f541279fc806ce39ecf1a53293d5c57c22ffd4/e2e/d625f10b-3b86-47c4-a5ca-ee8f4c1a7543.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/032a00e7b8b546b20e5ae53fa8d86a473852b56c/e2e/d625f10b-3b86-47c4-a5ca-ee8f4c1a7543.md."

# ---------------------------------------------------------------------------
# Overview sheet: rows 4 (75020932) and 5 (d625f10b)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E4").Value = $readyForHandoff
$wsOverview.Range("F4").Value = $readyForHandoff
$wsOverview.Range("G4").Value = "2016-08-23 12:25:19"

$wsOverview.Range("E5").Value = $readyForHandoff
$wsOverview.Range("F5").Value = $readyForHandoff
$wsOverview.Range("G5").Value = "2016-08-23 12:25:19"

# ---------------------------------------------------------------------------
# zh-cn sheet: rows 4 (75020932) and 5 (d625f10b)
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C4").Value = $readyForHandoff
$wsZhCn.Range("H4").Value = "2016-08-23 12:25:11"
$wsZhCn.Range("P4").Value = $msg75020932

$wsZhCn.Range("C5").Value = $readyForHandoff
$wsZhCn.Range("H5").Value = "2016-08-23 12:25:11"
$wsZhCn.Range("P5").Value = $msgD625f10b

# 39.2 (Excel's "characters" column-width unit) round-trips to a stored
# OOXML width of exactly 40, matching the target width of column P.
$wsZhCn.Columns.Item(16).ColumnWidth = 39.2

# ---------------------------------------------------------------------------
# de-de sheet: rows 4 (75020932) and 5 (d625f10b)
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C4").Value = $readyForHandoff
$wsDeDe.Range("H4").Value = "2016-08-23 12:25:19"
$wsDeDe.Range("P4").Value = $msg75020932

$wsDeDe.Range("C5").Value = $readyForHandoff
$wsDeDe.Range("H5").Value = "2016-08-23 12:25:19"
$wsDeDe.Range("P5").Value = $msgD625f10b

# 39.2 (Excel's "characters" column-width unit) round-trips to a stored
# OOXML width of exactly 40, matching the target width of column P.
$wsDeDe.Columns.Item(16).ColumnWidth = 39.2
